# fix typos in SMPTE boilerplate
#
# 1) Give the smpte-indent / smpte-note / smpte-example example paragraphs
#    an explicit (direct-formatting) indent that overrides the indent
#    coming from their paragraph style.
# 2) Apply a 90% horizontal character scaling to the Consolas-based
#    syntax-highlighting character styles (the "...Tok" styles) plus
#    "Verbatim Char", used by the SourceCode block.

$d = $word.ActiveDocument

foreach ($p in $d.Paragraphs) {
    $styleName = $p.Range.ParagraphFormat.Style.NameLocal

    if ($styleName -eq "smpte-indent") {
        $p.Range.ParagraphFormat.LeftIndent = 28.35
    }
    elseif ($styleName -eq "smpte-note") {
        # Only the worked example paragraph (not the trailing blank
        # smpte-note paragraph further down) gets the direct indent.
        if ($p.Range.Text.Trim().Length -gt 0) {
            $p.Range.ParagraphFormat.LeftIndent = 28.35
            $p.Range.ParagraphFormat.FirstLineIndent = -28.9
        }
    }
    elseif ($styleName -eq "smpte-example") {
        $p.Range.ParagraphFormat.LeftIndent = 56.7
        $p.Range.ParagraphFormat.FirstLineIndent = -28.9
    }
}

$tokStyleNames = @(
    "KeywordTok", "DataTypeTok", "DecValTok", "BaseNTok", "FloatTok",
    "ConstantTok", "CharTok", "SpecialCharTok", "StringTok",
    "VerbatimStringTok", "SpecialStringTok", "ImportTok", "CommentTok",
    "DocumentationTok", "AnnotationTok", "CommentVarTok", "OtherTok",
    "FunctionTok", "VariableTok", "ControlFlowTok", "OperatorTok",
    "BuiltInTok", "ExtensionTok", "PreprocessorTok", "AttributeTok",
    "RegionMarkerTok", "InformationTok", "WarningTok", "AlertTok",
    "ErrorTok", "NormalTok", "Verbatim Char"
)

foreach ($name in $tokStyleNames) {
    $s = $d.Styles.Item($name)
    $s.Font.Scaling = 90
}
